$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date (column C) for rows 2-7 from 2023-09-01 (45170) to 2023-09-05 (45174)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45174
}
